# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" sheet between "2021-Q3" and "总计", built on the
#    layout of the "2021-Q3" fund-holding sheet (same headers/styles) but
#    with the fund numbers updated for the new quarter.
# 2. Insert a new leading row into the "总计" (totals) sheet recording the
#    2022-Q1 holding summary, pushing the existing rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a LITERAL TEXT value into a cell without Excel's automatic
# "numeric-looking string -> number" coercion (and without leaving a
# quote-prefix marker behind). We do this by writing a text-producing
# formula into a scratch cell, copying its calculated value (values-only)
# onto the destination, then clearing the scratch cell.
# ---------------------------------------------------------------------------
function Set-TextValue($ws, $range, [string]$text) {
    $helper = $ws.Range("ZZ100")
    $escaped = $text.Replace('"', '""')
    $helper.Formula = '="' + $escaped + '"'
    $helper.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null   # xlPasteValues
    $helper.Clear() | Out-Null
    $ws.Application.CutCopyMode = $false
}

# --- 1. Insert the new "2022-Q1" sheet right after "2021-Q3" ---------------
$q3 = $wb.Worksheets.Item("2021-Q3")
$newQ1 = $wb.Worksheets.Add($null, $q3)
$newQ1.Name = "2022-Q1"

# Base the new sheet on a full copy (values + formats) of "2021-Q3", which
# already has the right header/style layout for a quarter holdings sheet.
$q3.Range("A1:H2").Copy() | Out-Null
$newQ1.Range("A1:H2").PasteSpecial(-4104) | Out-Null   # xlPasteAll -> values
$q3.Range("A1:H2").Copy() | Out-Null
$newQ1.Range("A1:H2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> styles
$excel.CutCopyMode = $false

# Header: column D is relabelled from "基金金额" to "基金规模"
$newQ1.Range("D1").Value = "基金规模"

# Data row: update the figures that changed for 2022-Q1
# (B2 "968029", C2 fund name and H2 "8" stay the same as 2021-Q3)
Set-TextValue $newQ1 $newQ1.Range("D2") "25.09"
Set-TextValue $newQ1 $newQ1.Range("E2") "97.94"
Set-TextValue $newQ1 $newQ1.Range("F2") "2.95"
Set-TextValue $newQ1 $newQ1.Range("G2") "0.7402"

# --- 2. Record the new quarter in the "总计" (totals) sheet ----------------
$total = $wb.Worksheets.Item("总计")
$total.Rows("2:2").Insert() | Out-Null

# The row-insert can drag stray formatting into the new blank row; make sure
# B2:D2 start from a clean (unstyled) slate like the rest of the data rows.
$total.Range("B2:D2").ClearFormats() | Out-Null

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.74

# A2 should carry the same index style as the other index cells (A3, A4 ...)
$total.Range("A3").Copy() | Out-Null
$total.Range("A2").PasteSpecial(-4122) | Out-Null      # xlPasteFormats -> styles
$excel.CutCopyMode = $false

# Renumber the index column for the rows that shifted down
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2

Write-Output "Added 2022-Q1 sheet and updated totals"
